$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B (epoch accuracy) values that changed between runs
$bUpdates = @{
    3 = 0.828125
    4 = 0.75
    5 = 0.671875
    6 = 0.59375
    7 = 0.609375
    8 = 0.546875
    9 = 0.546875
    10 = 0.53125
    11 = 0.59375
    12 = 0.59375
    13 = 0.546875
    14 = 0.5625
    16 = 0.59375
    17 = 0.59375
    18 = 0.59375
    19 = 0.578125
    20 = 0.578125
    21 = 0.59375
    22 = 0.578125
    23 = 0.5625
    24 = 0.5625
    25 = 0.546875
    26 = 0.546875
    27 = 0.546875
    28 = 0.546875
    29 = 0.546875
    30 = 0.546875
    31 = 0.546875
    32 = 0.546875
    33 = 0.546875
    34 = 0.546875
    35 = 0.546875
    36 = 0.546875
    37 = 0.546875
    38 = 0.546875
    39 = 0.546875
    40 = 0.546875
    41 = 0.546875
    42 = 0.546875
    43 = 0.546875
    44 = 0.546875
    58 = 0.5625
    59 = 0.5625
    60 = 0.5625
    61 = 0.5625
    62 = 0.5625
    63 = 0.5625
    64 = 0.5625
    65 = 0.5625
    66 = 0.5625
    67 = 0.5625
    68 = 0.5625
    69 = 0.5625
    70 = 0.5625
    71 = 0.5625
    72 = 0.5625
    73 = 0.5625
    74 = 0.5625
    75 = 0.5625
    76 = 0.5625
    77 = 0.5625
    78 = 0.5625
    79 = 0.5625
    80 = 0.5625
    81 = 0.5625
    82 = 0.5625
    83 = 0.5625
    84 = 0.5625
    85 = 0.5625
    86 = 0.5625
    87 = 0.5625
    88 = 0.5625
    89 = 0.5625
    90 = 0.5625
    91 = 0.5625
    92 = 0.5625
    93 = 0.5625
    94 = 0.5625
    95 = 0.5625
    96 = 0.5625
    97 = 0.5625
    98 = 0.5625
    99 = 0.5625
    100 = 0.5625
    101 = 0.5625
    102 = 0.5625
    103 = 0.546875
    104 = 0.53125
    105 = 0.5625
    107 = 0.5625
    108 = 0.5625
    110 = 0.640625
    111 = 0.578125
    112 = 0.5
    113 = 0.703125
    114 = 0.5625
    115 = 0.609375
    116 = 0.6176470588235294
}
foreach ($row in $bUpdates.Keys) {
    $ws.Cells.Item($row, 2).Value = $bUpdates[$row]
}

# Column A rows 102-116 hold the repr() of a DisplayOutputs object; only the
# in-memory address changed between the two interpreter runs
for ($r = 102; $r -le 116; $r++) {
    $ws.Cells.Item($r, 1).Value = "<__main__.DisplayOutputs object at 0x7f21e36af760>"
}

# Reselect: user selected the whole sheet (Ctrl+A), which Excel records as A1:XFD1048576
$null = $ws.Cells.Select()
